{"js": "// Replace the two-digit multiplication problems/answers with the new\n// values from the commit, matching each old \"A\u00d7B=C\" text exactly and\n// swapping in the new \"A\u00d7B=C\" text.\nconst replacements = [\n  [\"45\u00d775=3375\", \"44\u00d795=4180\"],\n  [\"31\u00d774=2294\", \"75\u00d733=2475\"],\n  [\"13\u00d766=858\", \"78\u00d781=6318\"],\n  [\"92\u00d742=3864\", \"39\u00d729=1131\"],\n  [\"22\u00d714=308\", \"27\u00d714=378\"],\n  [\"91\u00d764=5824\", \"18\u00d728=504\"],\n  [\"45\u00d785=3825\", \"62\u00d795=5890\"],\n  [\"84\u00d743=3612\", \"46\u00d735=1610\"],\n  [\"14\u00d756=784\", \"65\u00d715=975\"],\n  [\"47\u00d741=1927\", \"15\u00d736=540\"],\n  [\"37\u00d793=3441\", \"32\u00d717=544\"],\n  [\"83\u00d758=4814\", \"52\u00d778=4056\"],\n  [\"76\u00d769=5244\", \"37\u00d766=2442\"],\n  [\"61\u00d755=3355\", \"11\u00d776=836\"],\n  [\"79\u00d749=3871\", \"70\u00d750=3500\"],\n  [\"97\u00d753=5141\", \"82\u00d725=2050\"],\n  [\"48\u00d723=1104\", \"74\u00d734=2516\"],\n  [\"37\u00d773=2701\", \"15\u00d720=300\"],\n  [\"55\u00d788=4840\", \"58\u00d798=5684\"],\n  [\"42\u00d749=2058\", \"81\u00d752=4212\"],\n  [\"17\u00d775=1275\", \"84\u00d771=5964\"],\n  [\"53\u00d740=2120\", \"77\u00d762=4774\"],\n  [\"31\u00d725=775\", \"76\u00d761=4636\"],\n  [\"23\u00d792=2116\", \"79\u00d719=1501\"],\n  [\"78\u00d713=1014\", \"73\u00d784=6132\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit multiplication problems/answers with the new\n# values from the commit, matching each old \"A\u00d7B=C\" text exactly and\n# swapping in the new \"A\u00d7B=C\" text.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"45\u00d775=3375\"; New = \"44\u00d795=4180\" },\n    @{ Old = \"31\u00d774=2294\"; New = \"75\u00d733=2475\" },\n    @{ Old = \"13\u00d766=858\";  New = \"78\u00d781=6318\" },\n    @{ Old = \"92\u00d742=3864\"; New = \"39\u00d729=1131\" },\n    @{ Old = \"22\u00d714=308\";  New = \"27\u00d714=378\" },\n    @{ Old = \"91\u00d764=5824\"; New = \"18\u00d728=504\" },\n    @{ Old = \"45\u00d785=3825\"; New = \"62\u00d795=5890\" },\n    @{ Old = \"84\u00d743=3612\"; New = \"46\u00d735=1610\" },\n    @{ Old = \"14\u00d756=784\";  New = \"65\u00d715=975\" },\n    @{ Old = \"47\u00d741=1927\"; New = \"15\u00d736=540\" },\n    @{ Old = \"37\u00d793=3441\"; New = \"32\u00d717=544\" },\n    @{ Old = \"83\u00d758=4814\"; New = \"52\u00d778=4056\" },\n    @{ Old = \"76\u00d769=5244\"; New = \"37\u00d766=2442\" },\n    @{ Old = \"61\u00d755=3355\"; New = \"11\u00d776=836\" },\n    @{ Old = \"79\u00d749=3871\"; New = \"70\u00d750=3500\" },\n    @{ Old = \"97\u00d753=5141\"; New = \"82\u00d725=2050\" },\n    @{ Old = \"48\u00d723=1104\"; New = \"74\u00d734=2516\" },\n    @{ Old = \"37\u00d773=2701\"; New = \"15\u00d720=300\" },\n    @{ Old = \"55\u00d788=4840\"; New = \"58\u00d798=5684\" },\n    @{ Old = \"42\u00d749=2058\"; New = \"81\u00d752=4212\" },\n    @{ Old = \"17\u00d775=1275\"; New = \"84\u00d771=5964\" },\n    @{ Old = \"53\u00d740=2120\"; New = \"77\u00d762=4774\" },\n    @{ Old = \"31\u00d725=775\";  New = \"76\u00d761=4636\" },\n    @{ Old = \"23\u00d792=2116\"; New = \"79\u00d719=1501\" },\n    @{ Old = \"78\u00d713=1014\"; New = \"73\u00d784=6132\" }\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n}\n"}
